$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ID -> 3, FilePath -> CloneScene, SceneName -> clone
$ws.Range("B2").Value = "3"
$ws.Range("A2").Value = "../../NFDataCfg/Ini/NFZoneServer/Scene/CloneScene/"
$ws.Range("F2").Value = "clone"

# Row 3: SceneName -> newscene (ID stays 1, FilePath stays PioneerNoob)
$ws.Range("F3").Value = "newscene"

# Row 4: SceneName -> newscene (ID stays 2, FilePath stays RebellerNoob)
$ws.Range("F4").Value = "newscene"

# Update the selected cell to match the new active selection in the diff
$ws.Range("H8").Select()
